# Updates FFXIV leve-profit calculations across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per scheduled market-price data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 226.93333
$ws.Range("I33").Value = 231.76923
$ws.Range("J33").Value = 195.5
$ws.Range("K33").Value = 231.76923
$ws.Range("L33").Value = 195.5
$ws.Range("M33").Value = -2.769229999999993
$ws.Range("N33").Value = -653.5

$ws.Range("H70").Value = 2000.2858
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 2000.3334
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 6001.0002
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -6541.0002

$ws.Range("H73").Value = 2000.2858
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 2000.3334
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 6001.0002
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -7873.0002

$ws.Range("H103").Value = 910.05
$ws.Range("I103").Value = 478.57144
$ws.Range("J103").Value = 1142.3846
$ws.Range("K103").Value = 1435.71432
$ws.Range("L103").Value = 3427.1538
$ws.Range("M103").Value = -849.71432
$ws.Range("N103").Value = -4599.1538

$ws.Range("H125").Value = 808.125
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 910.8333
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 8197.4997
$ws.Range("M125").Value = -2040
$ws.Range("N125").Value = -13117.4997

$ws.Range("H129").Value = 927.59015
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 971.6316
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 2914.8948
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -12914.8948

$ws.Range("H138").Value = 1681.09
$ws.Range("I138").Value = 850.55554
$ws.Range("J138").Value = 1863.4025
$ws.Range("K138").Value = 2551.66662
$ws.Range("L138").Value = 5590.2075
$ws.Range("M138").Value = 2588.33338
$ws.Range("N138").Value = -15870.2075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4999.8
$ws.Range("I61").Value = 4029.05
$ws.Range("J61").Value = 6941.3
$ws.Range("K61").Value = 4029.05
$ws.Range("L61").Value = 6941.3
$ws.Range("M61").Value = -3817.05
$ws.Range("N61").Value = -7365.3

$ws.Range("H122").Value = 4883
$ws.Range("I122").Value = 4079.2222
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 12237.6666
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -9787.6666
$ws.Range("N122").Value = -30400

$ws.Range("H132").Value = 2600.8125
$ws.Range("I132").Value = 1628.4
$ws.Range("J132").Value = 3458.8235
$ws.Range("K132").Value = 4885.200000000001
$ws.Range("L132").Value = 10376.4705
$ws.Range("M132").Value = -2355.200000000001
$ws.Range("N132").Value = -15436.4705

$ws.Range("H136").Value = 4999.8
$ws.Range("I136").Value = 4029.05
$ws.Range("J136").Value = 6941.3
$ws.Range("K136").Value = 12087.15
$ws.Range("L136").Value = 20823.9
$ws.Range("M136").Value = -9537.150000000001
$ws.Range("N136").Value = -25923.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 40271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 40271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41643

$ws.Range("H66").Value = 40271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 40271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127677

$ws.Range("H132").Value = 49141
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 49141
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 49141
$ws.Range("N132").Value = -59261

$ws.Range("H134").Value = 6742.2593
$ws.Range("I134").Value = 10498.134
$ws.Range("J134").Value = 2047.4166
$ws.Range("K134").Value = 31494.402
$ws.Range("L134").Value = 6142.2498
$ws.Range("M134").Value = -28959.402
$ws.Range("N134").Value = -11212.2498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H31").Value = 3263.8965
$ws.Range("I31").Value = 2385.6667
$ws.Range("J31").Value = 4701
$ws.Range("K31").Value = 2385.6667
$ws.Range("L31").Value = 4701
$ws.Range("M31").Value = -2090.6667
$ws.Range("N31").Value = -5291

$ws.Range("H34").Value = 3263.8965
$ws.Range("I34").Value = 2385.6667
$ws.Range("J34").Value = 4701
$ws.Range("K34").Value = 2385.6667
$ws.Range("L34").Value = 4701
$ws.Range("M34").Value = -2183.6667
$ws.Range("N34").Value = -5105

$ws.Range("H39").Value = 10279.2
$ws.Range("I39").Value = 4149
$ws.Range("J39").Value = 34800
$ws.Range("K39").Value = 4149
$ws.Range("L39").Value = 34800
$ws.Range("M39").Value = -35582
$ws.Range("N39").Value = -35582

$ws.Range("H49").Value = 10279.2
$ws.Range("I49").Value = 4149
$ws.Range("J49").Value = 34800
$ws.Range("K49").Value = 4149
$ws.Range("L49").Value = 34800
$ws.Range("M49").Value = -35164
$ws.Range("N49").Value = -35164

$ws.Range("H99").Value = 1697.3334
$ws.Range("I99").Value = 1600
$ws.Range("J99").Value = 1746
$ws.Range("K99").Value = 1600
$ws.Range("L99").Value = 1746
$ws.Range("M99").Value = -102
$ws.Range("N99").Value = -4742

$ws.Range("H126").Value = 1697.3334
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 1746
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 5238
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -10178

$ws.Range("H141").Value = 42352.4
$ws.Range("I141").Value = 50296
$ws.Range("J141").Value = 41934.316
$ws.Range("K141").Value = 50296
$ws.Range("L141").Value = 41934.316
$ws.Range("M141").Value = -45116
$ws.Range("N141").Value = -52294.316

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3790739
$ws.Range("I5").Value = 403.0303
$ws.Range("J5").Value = 15161747
$ws.Range("K5").Value = 1209.0909
$ws.Range("L5").Value = 45485241
$ws.Range("M5").Value = -1097.0909
$ws.Range("N5").Value = -45485465

$ws.Range("H41").Value = 150.25
$ws.Range("I41").Value = 150.25
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 450.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -112.75
$ws.Range("N41").ClearContents()

$ws.Range("H43").Value = 5999
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5999
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 17997
$ws.Range("N43").Value = -18225

$ws.Range("H135").Value = 3790739
$ws.Range("I135").Value = 403.0303
$ws.Range("J135").Value = 15161747
$ws.Range("K135").Value = 3627.2727
$ws.Range("L135").Value = 136455723
$ws.Range("M135").Value = -1092.2727
$ws.Range("N135").Value = -136460793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3764
$ws.Range("I102").Value = 2939.8333
$ws.Range("J102").Value = 7720
$ws.Range("K102").Value = 2939.8333
$ws.Range("L102").Value = 7720
$ws.Range("M102").Value = -1317.8333
$ws.Range("N102").Value = -10964

$ws.Range("H132").Value = 1357.0416
$ws.Range("I132").Value = 435.75
$ws.Range("J132").Value = 2278.3333
$ws.Range("K132").Value = 1307.25
$ws.Range("L132").Value = 6834.999899999999
$ws.Range("M132").Value = 1222.75
$ws.Range("N132").Value = -11894.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 387.77777
$ws.Range("I22").Value = 247.5
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 247.5
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 47.5
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 387.77777
$ws.Range("I27").Value = 247.5
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 247.5
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -140.5
$ws.Range("N27").Value = -714

$ws.Range("H46").Value = 901.5
$ws.Range("I46").Value = 801
$ws.Range("K46").Value = 801
$ws.Range("M46").Value = -613

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H132").Value = 3647.7742
$ws.Range("I132").Value = 3180.1428
$ws.Range("J132").Value = 4032.8823
$ws.Range("K132").Value = 9540.428400000001
$ws.Range("L132").Value = 12098.6469
$ws.Range("M132").Value = -7010.428400000001
$ws.Range("N132").Value = -17158.6469

$ws.Range("H136").Value = 4187.396
$ws.Range("I136").Value = 2989.4211
$ws.Range("J136").Value = 7222.2666
$ws.Range("K136").Value = 8968.263300000001
$ws.Range("L136").Value = 21666.7998
$ws.Range("M136").Value = -6418.263300000001
$ws.Range("N136").Value = -26766.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3908.3845
$ws.Range("I62").Value = 3971.4285
$ws.Range("J62").Value = 3834.8333
$ws.Range("K62").Value = 3971.4285
$ws.Range("L62").Value = 3834.8333
$ws.Range("M62").Value = -3347.4285
$ws.Range("N62").Value = -5082.8333

$ws.Range("H65").Value = 3908.3845
$ws.Range("I65").Value = 3971.4285
$ws.Range("J65").Value = 3834.8333
$ws.Range("K65").Value = 19857.1425
$ws.Range("L65").Value = 19174.1665
$ws.Range("M65").Value = -16737.1425
$ws.Range("N65").Value = -25414.1665

$ws.Range("H107").Value = 5312
$ws.Range("I107").Value = 1933.3334
$ws.Range("J107").Value = 7001.3335
$ws.Range("K107").Value = 5800.0002
$ws.Range("L107").Value = 21004.0005
$ws.Range("M107").Value = -3880.0002
$ws.Range("N107").Value = -24844.0005

$ws.Range("H122").Value = 2274.7144
$ws.Range("I122").Value = 2058.5417
$ws.Range("J122").Value = 2746.3635
$ws.Range("K122").Value = 6175.625100000001
$ws.Range("L122").Value = 8239.0905
$ws.Range("M122").Value = -3725.625100000001
$ws.Range("N122").Value = -13139.0905

$ws.Range("H132").Value = 1622.1951
$ws.Range("I132").Value = 842
$ws.Range("J132").Value = 2619.111
$ws.Range("K132").Value = 2526
$ws.Range("L132").Value = 7857.333
$ws.Range("M132").Value = 4
$ws.Range("N132").Value = -12917.333

$ws.Range("H136").Value = 4028.2576
$ws.Range("I136").Value = 1783.65
$ws.Range("J136").Value = 7481.5
$ws.Range("K136").Value = 5350.950000000001
$ws.Range("L136").Value = 22444.5
$ws.Range("M136").Value = -2800.950000000001
$ws.Range("N136").Value = -27544.5

$ws.Range("H137").Value = 48589.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 48589.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 48589.5
$ws.Range("N137").Value = -58789.5
